$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value could be auto-interpreted as a number by Excel;
# force them to remain plain text so they round-trip as inline/shared strings.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '67.346.40'
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('D3').Value = '3.513.58'
$ws.Range('E3').Value = '  -1.28%  '
$ws.Range('D5').Value = '610.13'
$ws.Range('E5').Value = '  -0.90%  '
$ws.Range('D6').Value = '150.31'
$ws.Range('E6').Value = '  -2.33%  '
$ws.Range('D7').Value = '3.512.75'
$ws.Range('E7').Value = '  -1.20%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '0.481'
$ws.Range('E9').Value = '  -1.14%  '
$ws.Range('E10').Value = '  -1.30%  '
$ws.Range('E11').Value = '  +2.00%  '
$ws.Range('E12').Value = '  -1.47%  '
$ws.Range('E13').Value = '  -2.11%  '
$ws.Range('D14').Value = '4.111.79'
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('D15').Value = '31.77'
$ws.Range('E15').Value = '  -0.83%  '
$ws.Range('D16').Value = '3.517.76'
$ws.Range('E16').Value = '  -1.05%  '
$ws.Range('D17').Value = '67.368.81'
$ws.Range('E17').Value = '  -0.70%  '
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('E20').Value = '  -2.66%  '
$ws.Range('D21').Value = '442.86'
$ws.Range('E21').Value = '  -2.53%  '
$ws.Range('E22').Value = '  -4.14%  '
$ws.Range('D23').Value = '0.624'
$ws.Range('E23').Value = '  -3.31%  '
$ws.Range('D24').Value = '77.21'
$ws.Range('E24').Value = '  -0.50%  '
$ws.Range('E25').Value = '  +10.35%  '
$ws.Range('D26').Value = '3.657.82'
$ws.Range('E26').Value = '  -1.27%  '
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('D28').Value = '10.24'
$ws.Range('E28').Value = '  -3.67%  '
$ws.Range('D29').Value = '8.36'
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('E30').Value = '  -2.29%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').Value = '1.54'
$ws.Range('E31').Value = '  -4.73%  '
$ws.Range('B32').Value = 'Binance-PegBSC-USD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('E33').Value = '  +3.52%  '
$ws.Range('D34').Value = '25.75'
$ws.Range('D35').Value = '6.14'
$ws.Range('E35').Value = '  -1.20%  '
$ws.Range('D36').Value = '3.509.83'
$ws.Range('E36').Value = '  -1.40%  '
$ws.Range('E37').Value = '  -3.84%  '
$ws.Range('D38').Value = '7.99'
$ws.Range('E38').Value = '  -0.98%  '
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').Value = '178.06'
$ws.Range('E40').Value = '  +0.85%  '
$ws.Range('D42').Value = '2.16'
$ws.Range('E42').Value = '  +3.44%  '
$ws.Range('D43').Value = '0.0870'
$ws.Range('E43').Value = '  -1.65%  '
$ws.Range('E44').Value = '  -3.17%  '
$ws.Range('D45').Value = '0.880'
$ws.Range('E45').Value = '  -1.49%  '
$ws.Range('D46').Value = '45.57'
$ws.Range('E46').Value = '  -1.48%  '
$ws.Range('D47').Value = '27.46'
$ws.Range('E47').Value = '  -4.28%  '
$ws.Range('E48').Value = '  +5.11%  '
$ws.Range('D49').Value = '2.60'
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('E50').Value = '  -1.70%  '
$ws.Range('D51').Value = '0.996'
$ws.Range('E51').Value = '  -1.35%  '
